$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("N1").Value = "Demand Charges Block 1 Limit (MW)"
$ws.Range("P1").Value = "Demand Charges Block 2 Limit (MW)"
$ws.Range("R1").Value = "Demand Charges Block 3 Limit (MW)"
$ws.Range("T1").Value = "Demand Charges Block 4 Limit (MW)"
